# Update distribution stats: drop the two intermediate iteration rows (2 and 3)
# and refresh the final row (row 4) with the recalculated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 2 and 3 data (intermediate samples no longer kept)
$ws.Rows("2:3").ClearContents()

# Refresh row 4 with the updated gamma-distribution summary values
$ws.Range("A4").Value = 0.3300548414179534
$ws.Range("B4").Value = 0.438391105008159
$ws.Range("C4").Value = 0.2679979988735469
$ws.Range("D4").Value = 0.3937391041631541
$ws.Range("E4").Value = 0.2828853923370157
$ws.Range("F4").Value = 7.007322609896816
$ws.Range("G4").Value = 10.59244524579269
$ws.Range("H4").Value = 4.943887736082482
$ws.Range("I4").Value = 9.156663608614821
$ws.Range("J4").Value = 5.454697200150446
$ws.Range("K4").Value = 2.312828591780329
$ws.Range("L4").Value = 3.776206511706212
$ws.Range("M4").Value = 1.504166808897832
$ws.Range("N4").Value = 3.157386487208637
$ws.Range("O4").Value = 1.705213153564048
$ws.Range("P4").Value = 65.09828
$ws.Range("Q4").Value = 90.86152396822476
$ws.Range("R4").Value = 47.47755553305626
$ws.Range("S4").Value = 74.3787818779372
$ws.Range("T4").Value = 54.8664087005132
$ws.Range("U4").Value = 0.2014188176078627
$ws.Range("V4").Value = 0.3170577056105753
$ws.Range("W4").Value = 0.09951348041250996
$ws.Range("X4").Value = 0.2601819905801135
$ws.Range("Y4").Value = 0.1428659134014484
$ws.Range("Z4").Value = 0.8480062911440618
$ws.Range("AA4").Value = 0.9544705393367487
$ws.Range("AB4").Value = 0.6981762718042079
$ws.Range("AC4").Value = 0.9096999253447704
$ws.Range("AD4").Value = 0.7729972162013145
